$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove existing hyperlinks; we will re-add them in final row order
$ws.Hyperlinks.Delete()

# Row 2: 【急募】掲示板サイト(爆サイ)でAIによる自然な会話で書き込みを埋めていけるソフ
$ws.Cells.Item(2,1).Value = '2025-11-28 18:25:20'
$ws.Cells.Item(2,2).Value = '【急募】掲示板サイト(爆サイ)でAIによる自然な会話で書き込みを埋めていけるソフト開発者募集'
$ws.Cells.Item(2,3).Value = 'システム開発'
$ws.Cells.Item(2,4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(2,5).Value = '期限情報なし'
$ws.Cells.Item(2,6).Value = 'https://www.lancers.jp/work/detail/5443464'
$ws.Cells.Item(2,7).Value = 378
$ws.Cells.Item(2,8).Value = '🔥AI,Ai ◆開発 ◇サイト'
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5443464')

# Row 3: エッジAIカメラによる人流計測のPoC用プログラム開発
$ws.Cells.Item(3,1).Value = '2025-11-28 18:25:20'
$ws.Cells.Item(3,2).Value = 'エッジAIカメラによる人流計測のPoC用プログラム開発'
$ws.Cells.Item(3,3).Value = 'システム開発'
$ws.Cells.Item(3,4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(3,5).Value = '期限情報なし'
$ws.Cells.Item(3,6).Value = 'https://www.lancers.jp/work/detail/5443336'
$ws.Cells.Item(3,7).Value = 368
$ws.Cells.Item(3,8).Value = '🔥AI,Ai ◆開発'
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5443336')

# Row 4: 【急募】フロントエンド開発者募集!React/TypeScriptでのシステム構
$ws.Cells.Item(4,1).Value = '2025-11-28 18:25:20'
$ws.Cells.Item(4,2).Value = '【急募】フロントエンド開発者募集!React/TypeScriptでのシステム構築'
$ws.Cells.Item(4,3).Value = 'システム開発'
$ws.Cells.Item(4,4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(4,5).Value = '期限情報なし'
$ws.Cells.Item(4,6).Value = 'https://www.lancers.jp/work/detail/5443491'
$ws.Cells.Item(4,7).Value = 323
$ws.Cells.Item(4,8).Value = '🔥React,TypeScript ◆開発'
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5443491')

# Row 5: 【自動化】エクセルデータ転記作業の効率化依頼
$ws.Cells.Item(5,1).Value = '2025-11-28 18:25:20'
$ws.Cells.Item(5,2).Value = '【自動化】エクセルデータ転記作業の効率化依頼'
$ws.Cells.Item(5,3).Value = 'システム開発'
$ws.Cells.Item(5,4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(5,5).Value = '期限情報なし'
$ws.Cells.Item(5,6).Value = 'https://www.lancers.jp/work/detail/5442971'
$ws.Cells.Item(5,7).Value = 145
$ws.Cells.Item(5,8).Value = '◆効率化,自動化'
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5442971')

# Row 6: MT5アラートツールの制作
$ws.Cells.Item(6,1).Value = '2025-11-28 18:25:20'
$ws.Cells.Item(6,2).Value = 'MT5アラートツールの制作'
$ws.Cells.Item(6,3).Value = 'システム開発'
$ws.Cells.Item(6,4).Value = '10,000 円 ~ 20,000 円 / 募集期間 3 日、取引期間 0 日'
$ws.Cells.Item(6,5).Value = '期限情報なし'
$ws.Cells.Item(6,6).Value = 'https://www.lancers.jp/work/detail/5443470'
$ws.Cells.Item(6,7).Value = 65
$ws.Cells.Item(6,8).Value = '◆ツール'
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5443470')

# Row 7: 初回 2026年1月創業 コンサル会社のバックオフィス業務フロー設計・マニュアル
$ws.Cells.Item(7,1).Value = '2025-11-28 18:25:20'
$ws.Cells.Item(7,2).Value = '初回 2026年1月創業 コンサル会社のバックオフィス業務フロー設計・マニュアル化、IT導入 一括見積依頼'
$ws.Cells.Item(7,3).Value = 'システム開発'
$ws.Cells.Item(7,4).Value = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Cells.Item(7,5).Value = '期限情報なし'
$ws.Cells.Item(7,6).Value = 'https://www.lancers.jp/work/detail/5442904'
$ws.Cells.Item(7,7).Value = 55
$ws.Cells.Item(7,8).Value = '◆コンサル'
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5442904')

# Row 8: 【急募】PGエンジニア募集!構成管理・マスタメンテ業務
$ws.Cells.Item(8,1).Value = '2025-11-28 18:25:20'
$ws.Cells.Item(8,2).Value = '【急募】PGエンジニア募集!構成管理・マスタメンテ業務'
$ws.Cells.Item(8,3).Value = 'システム開発'
$ws.Cells.Item(8,4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(8,5).Value = '期限情報なし'
$ws.Cells.Item(8,6).Value = 'https://www.lancers.jp/work/detail/5443303'
$ws.Cells.Item(8,7).Value = 45
$ws.Cells.Item(8,8).Value = '◇管理'
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5443303')

# Row 9: 【急募】古いPHPとPerlプログラムのアップデート依頼
$ws.Cells.Item(9,1).Value = '2025-11-28 18:25:20'
$ws.Cells.Item(9,2).Value = '【急募】古いPHPとPerlプログラムのアップデート依頼'
$ws.Cells.Item(9,3).Value = 'システム開発'
$ws.Cells.Item(9,4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(9,5).Value = '期限情報なし'
$ws.Cells.Item(9,6).Value = 'https://www.lancers.jp/work/detail/5440861'
$ws.Cells.Item(9,7).Value = 33
$ws.Cells.Item(9,8).Value = '○PHP'
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5440861')

# Row 10: 【急募】PSE認証代行をお手伝いしてくれる方募集!
$ws.Cells.Item(10,1).Value = '2025-11-28 18:25:20'
$ws.Cells.Item(10,2).Value = '【急募】PSE認証代行をお手伝いしてくれる方募集!'
$ws.Cells.Item(10,3).Value = 'システム開発'
$ws.Cells.Item(10,4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(10,5).Value = '期限情報なし'
$ws.Cells.Item(10,6).Value = 'https://www.lancers.jp/work/detail/5443188'
$ws.Cells.Item(10,7).Value = 10
$ws.Cells.Item(10,8).Value = ""
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5443188')

# Column width adjustments (target raw OOXML width = ColumnWidth + 5/6)
$ws.Columns("D").ColumnWidth = 41 - 5/6
$ws.Columns("H").ColumnWidth = 23 - 5/6

